$d = $word.ActiveDocument

# =====================================================================
# 1. First paragraph: append "  (This is a change – Version for main
#    branch)" in red, split across three runs the same way the
#    reference edit does.
# =====================================================================
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$textRange = $d.Range($r.Start, $r.End - 1)
$textRange.InsertAfter("  ")

$p1 = $d.Paragraphs.Item(1)
$insertPoint = $p1.Range.End - 1
$seg1 = $d.Range($insertPoint, $insertPoint)
$seg1.InsertAfter("(This is a change – Ve")
$seg1.Font.Color = 255

$p1 = $d.Paragraphs.Item(1)
$insertPoint2 = $p1.Range.End - 1
$seg2 = $d.Range($insertPoint2, $insertPoint2)
$seg2.InsertAfter("rsion for main branch")
$seg2.Font.Color = 255

$p1 = $d.Paragraphs.Item(1)
$insertPoint3 = $p1.Range.End - 1
$seg3 = $d.Range($insertPoint3, $insertPoint3)
$seg3.InsertAfter(")")
$seg3.Font.Color = 255

# =====================================================================
# 2. Remove the trailing "ank God almighty, we are free at last."
#    paragraph (the final paragraph in the document).
# =====================================================================
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastP.Range.Delete()

# =====================================================================
# 3. Drop the now-unused custom/heading styles that were removed from
#    styles.xml.
# =====================================================================
$stylesToDelete = @(
    "Heading 2",
    "Heading 4",
    "apple-converted-space",
    "Hyperlink",
    "Heading 2 Char",
    "Heading 4 Char",
    "audio-tool",
    "subscribe",
    "subscribe-more-info",
    "generic-title",
    "podcast-tools__subscribe-links"
)
foreach ($name in $stylesToDelete) {
    $st = $d.Styles.Item($name)
    $st.Delete()
}

Write-Output "done"
